# FFL_Data.xlsx edit: "Week 10" results added, "London Bridge is Down" team
# renamed to "Dulcich de Leche" (so every historical standings row that used
# to reference "London Bridge is Down" now reads "Dulcich de Leche"), and the
# Week 9 image references were cleared (the Image column data moved to the
# newly appended Week 10 rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Append the "Week 10" results block (rows 110-121) first.
# ---------------------------------------------------------------------------
$week10 = @(
  @("Kauaireek Hill",                   "Week 10", 76.7,   125.24, "Team Icons/kauaireek-modified.png"),
  @("Chasing dank Herb",                "Week 10", 118.82, 79.8,   "Team Icons/chasing-modified.png"),
  @("Ju Ju Smith Poopster",             "Week 10", 62.54,  120.06, "Team Icons/juju-modified.png"),
  @("Bye Breece See You in ValHalla",   "Week 10", 119.2,  127.48, "Team Icons/breece-modified.png"),
  @("Cooking with Gas",                 "Week 10", 120.06, 62.54,  "Team Icons/cooking-modified.png"),
  @("Dulcich de Leche",                 "Week 10", 96.3,   89.1,   "Team Icons/dulcich-modified.png"),
  @("Dillon Panthers",                  "Week 10", 125.24, 76.7,   "Team Icons/dillon-modified.png"),
  @("Daemon and the Rightful Heirs",    "Week 10", 127.48, 118.17, "Team Icons/daemon-modified.png"),
  @("Krombopulos Michael Evans",        "Week 10", 87.24,  125.68, "Team Icons/krombopulos-modified.png"),
  @("Freier Freier Pants on Fire",      "Week 10", 89.1,   96.3,   "Team Icons/freier-modified.png"),
  @("Sir Trey Lancelot (IR)",           "Week 10", 125.68, 87.24,  "Team Icons/sir-modified.png"),
  @("Fantasy Football Champion 2022",   "Week 10", 79.8,   118.82, "Team Icons/fantasy-modified.png")
)

$r = 110
foreach ($row in $week10) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 2) Rename the team "London Bridge is Down" -> "Dulcich de Leche" everywhere
#    it appears in the existing Week 1-9 standings (column A).
# ---------------------------------------------------------------------------
$teamRange = $ws.Range("A1:A109")
$teamRange.Replace("London Bridge is Down", "Dulcich de Leche", 1, 1, $false, $false, $false, $false)

# ---------------------------------------------------------------------------
# 3) Clear the stale per-week team-icon image references that used to sit in
#    the Week 9 block (rows 98-109) - the Image column content now belongs
#    to the newly added Week 10 block instead.
# ---------------------------------------------------------------------------
$ws.Range("E98:E109").ClearContents()

# ---------------------------------------------------------------------------
# 4) Restore the view state: frozen header row, scrolled down to the new
#    Week 10 block, with the last-used cell selected.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("A98").Select()
$win.FreezePanes = $true
$ws.Range("E127").Select()
